$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) values for D, J, K, L, M, P in rows 2..10
# before we start overwriting them, since the edit shifts each row's data
# down into the next row and appends a brand new row 11 built from the
# former row 10. Use .Value2 (raw primitive) rather than .Value so the
# captured numbers/strings round-trip cleanly when written back.
$rows = 2..10
$orig = @{}
foreach ($r in $rows) {
    $rec = @{
        D = $ws.Cells.Item($r, 4).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        P = $ws.Cells.Item($r, 16).Value2
    }
    $orig[$r] = $rec
}

# Also remember row 10's "static" columns (the ones that don't shift) so we
# can build the new row 11 from them.
$a10 = $ws.Cells.Item(10, 1).Value2
$b10 = $ws.Cells.Item(10, 2).Value2
$c10 = $ws.Cells.Item(10, 3).Value2
$e10 = $ws.Cells.Item(10, 5).Value2
$f10 = $ws.Cells.Item(10, 6).Value2
$g10 = $ws.Cells.Item(10, 7).Value2
$h10 = $ws.Cells.Item(10, 8).Value2
$i10 = $ws.Cells.Item(10, 9).Value2
$n10 = $ws.Cells.Item(10, 14).Value2
$o10 = $ws.Cells.Item(10, 15).Value2
$q10 = $ws.Cells.Item(10, 17).Value2
$r10 = $ws.Cells.Item(10, 18).Value2

# Row 2 keeps its other values but gets a brand new date.
$ws.Cells.Item(2, 4).Value = 44487

# Rows 3..10 each inherit the D/J/K/L/M/P values that used to live one row
# above them (i.e. the whole block of weekly records shifts down by one).
for ($r = 3; $r -le 10; $r++) {
    $prev = $orig[$r - 1]
    $ws.Cells.Item($r, 4).Value = $prev.D
    $ws.Cells.Item($r, 10).Value = $prev.J
    $ws.Cells.Item($r, 11).Value = $prev.K
    $ws.Cells.Item($r, 12).Value = $prev.L
    $ws.Cells.Item($r, 13).Value = $prev.M
    $ws.Cells.Item($r, 16).Value = $prev.P
}

# A brand new row 11 is appended, carrying the same non-shifting columns as
# row 10 plus the D/J/K/L/M/P values that used to belong to row 10 (before
# the shift above moved new data into row 10).
$last = $orig[10]
$ws.Cells.Item(11, 1).Value = $a10
$ws.Cells.Item(11, 2).Value = $b10
$ws.Cells.Item(11, 3).Value = $c10
$ws.Cells.Item(11, 4).Value = $last.D
$ws.Cells.Item(11, 5).Value = $e10
$ws.Cells.Item(11, 6).Value = $f10
$ws.Cells.Item(11, 7).Value = $g10
$ws.Cells.Item(11, 8).Value = $h10
$ws.Cells.Item(11, 9).Value = $i10
$ws.Cells.Item(11, 10).Value = $last.J
$ws.Cells.Item(11, 11).Value = $last.K
$ws.Cells.Item(11, 12).Value = $last.L
$ws.Cells.Item(11, 13).Value = $last.M
$ws.Cells.Item(11, 14).Value = $n10
$ws.Cells.Item(11, 15).Value = $o10
$ws.Cells.Item(11, 16).Value = $last.P
$ws.Cells.Item(11, 17).Value = $q10
$ws.Cells.Item(11, 18).Value = $r10

# Match the D column's date number format/style used by the other rows.
$ws.Range("D11").NumberFormat = $ws.Range("D10").NumberFormat
